# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Update the "quantity" column (F) values on the NetDemand sheet to reflect
# the corrected calculation (uncon_planned_qty for future production days,
# produced qty kept for today / available inventory fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetDemand")

$ws.Range("F2").Value = -242
$ws.Range("F3").Value = -546
$ws.Range("F4").Value = -851
$ws.Range("F5").Value = -107
$ws.Range("F7").Value = -107
